$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - match style used by other header cells (bold, centered, bordered)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

$times = @(
    "2021-10-05 13:39:16.095654",
    "2021-10-05 13:39:16.095665",
    "2021-10-05 13:39:16.095668",
    "2021-10-05 13:39:16.095671",
    "2021-10-05 13:39:16.095674",
    "2021-10-05 13:39:16.095676",
    "2021-10-05 13:39:16.095679",
    "2021-10-05 13:39:16.095681",
    "2021-10-05 13:39:16.095684",
    "2021-10-05 13:39:16.095687",
    "2021-10-05 13:39:16.095689",
    "2021-10-05 13:39:16.095692",
    "2021-10-05 13:39:16.095694",
    "2021-10-05 13:39:16.095697",
    "2021-10-05 13:39:16.095699",
    "2021-10-05 13:39:16.095702",
    "2021-10-05 13:39:16.095705",
    "2021-10-05 13:39:16.095707",
    "2021-10-05 13:39:16.095710",
    "2021-10-05 13:39:16.095713",
    "2021-10-05 13:39:16.095715",
    "2021-10-05 13:39:16.095718",
    "2021-10-05 13:39:16.095720"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $times[$i]
}
